$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("G4").Value = 1.83
$ws.Range("H4").Value = 3.75
$ws.Range("I4").Value = 4.2
$ws.Range("K4").Value = 12
$ws.Range("U4").Value = 9.5
$ws.Range("AB4").Value = 13
$ws.Range("AD4").Value = 151
$ws.Range("AG4").Value = 13
$ws.Range("AH4").Value = 41

# Row 11
$ws.Range("G11").Value = 1.85
$ws.Range("H11").Value = 3.5
$ws.Range("I11").Value = 4.2
$ws.Range("J11").Value = 1.07
$ws.Range("K11").Value = 8.5
$ws.Range("R11").Value = 1.91
$ws.Range("S11").Value = 1.8
$ws.Range("T11").Value = 6.5
$ws.Range("X11").Value = 17
$ws.Range("AA11").Value = 6.5
$ws.Range("AC11").Value = 51
$ws.Range("AD11").Value = 351
$ws.Range("AE11").Value = 10
$ws.Range("AH11").Value = 41

# Row 33
$ws.Range("L33").Value = 1.11
$ws.Range("M33").Value = 6
$ws.Range("N33").Value = 1.4
$ws.Range("O33").Value = 2.75

# Row 34
$ws.Range("G34").Value = 1.33
$ws.Range("H34").Value = 5.5
$ws.Range("I34").Value = 6.25
$ws.Range("L34").Value = 1.11
$ws.Range("M34").Value = 6
$ws.Range("R34").Value = 1.73
$ws.Range("S34").Value = 2
$ws.Range("V34").Value = 9
$ws.Range("W34").Value = 9.5
$ws.Range("AA34").Value = 12
$ws.Range("AB34").Value = 19
$ws.Range("AC34").Value = 51
$ws.Range("AE34").Value = 23
$ws.Range("AF34").Value = 41
$ws.Range("AG34").Value = 21
$ws.Range("AH34").Value = 67

# Row 42
$ws.Range("J42").Value = 1.04
$ws.Range("K42").Value = 13
$ws.Range("L42").Value = 1.25
$ws.Range("M42").Value = 3.75
$ws.Range("N42").Value = 1.8
$ws.Range("O42").Value = 2

# Row 46
$ws.Range("G46").Value = 8
$ws.Range("H46").Value = 6
$ws.Range("I46").Value = 1.25
$ws.Range("K46").Value = 34
$ws.Range("V46").Value = 26
$ws.Range("W46").Value = 101
$ws.Range("X46").Value = 51
$ws.Range("AB46").Value = 17
$ws.Range("AH46").Value = 11
